# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# sheet with the latest scraped values (GitHub Actions data refresh).
# A couple of Price cells (D9, D34) are forced to Text format first so
# Excel doesn't silently drop a meaningful trailing zero by treating the
# string as a number (e.g. "0.550" -> 0.55).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.299.57"
$ws.Range("E2").Value = "  +2.50%  "
$ws.Range("D3").Value = "2.363.51"
$ws.Range("E3").Value = "  +0.74%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "0.676"
$ws.Range("E5").Value = "  +4.30%  "
$ws.Range("D6").Value = "238.68"
$ws.Range("E6").Value = "  +3.30%  "
$ws.Range("D7").Value = "73.52"
$ws.Range("E7").Value = "  +12.36%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.550"
$ws.Range("E9").Value = "  +20.83%  "
$ws.Range("E10").Value = "  +7.12%  "
$ws.Range("D11").Value = "29.53"
$ws.Range("E11").Value = "  +11.06%  "
$ws.Range("D12").Value = "0.107"
$ws.Range("E12").Value = "  +2.91%  "
$ws.Range("D13").Value = "2.716.49"
$ws.Range("E13").Value = "  +0.78%  "
$ws.Range("D14").Value = "16.85"
$ws.Range("E14").Value = "  +10.28%  "
$ws.Range("D15").Value = "6.73"
$ws.Range("E15").Value = "  +7.96%  "
$ws.Range("D16").Value = "0.907"
$ws.Range("E16").Value = "  +8.48%  "
$ws.Range("D17").Value = "2.358.67"
$ws.Range("E17").Value = "  +0.43%  "
$ws.Range("D18").Value = "44.209.01"
$ws.Range("E18").Value = "  +2.41%  "
$ws.Range("E19").Value = "  +5.57%  "
$ws.Range("D20").Value = "78.07"
$ws.Range("E20").Value = "  +6.18%  "
$ws.Range("E21").Value = "  +4.70%  "
$ws.Range("D22").Value = "256.03"
$ws.Range("E22").Value = "  +3.62%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("E24").Value = "  -3.90%  "
$ws.Range("D25").Value = "2.53"
$ws.Range("E25").Value = "  +3.64%  "
$ws.Range("D26").Value = "10.53"
$ws.Range("E26").Value = "  +7.09%  "
$ws.Range("D28").Value = "22.52"
$ws.Range("E28").Value = "  +1.26%  "
$ws.Range("E29").Value = "  +5.63%  "
$ws.Range("D30").Value = "173.17"
$ws.Range("E30").Value = "  -1.12%  "
$ws.Range("E31").Value = "  +3.71%  "
$ws.Range("D32").Value = "0.133"
$ws.Range("E32").Value = "  +5.83%  "
$ws.Range("E33").Value = "  +5.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0740"
$ws.Range("E34").Value = "  +7.85%  "
$ws.Range("D35").Value = "5.25"
$ws.Range("E35").Value = "  +5.49%  "
$ws.Range("D36").Value = "3.94"
$ws.Range("E36").Value = "  +10.41%  "
$ws.Range("E37").Value = "  -1.55%  "
$ws.Range("E38").Value = "  +0.81%  "
$ws.Range("D39").Value = "0.0272"
$ws.Range("E39").Value = "  +8.07%  "
$ws.Range("D40").Value = "19.68"
$ws.Range("E40").Value = "  +10.68%  "
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("D42").Value = "8.89"
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").Value = "1.25"
$ws.Range("E43").Value = "  +4.02%  "
$ws.Range("E44").Value = "  +4.24%  "
$ws.Range("E45").Value = "  +0.94%  "
$ws.Range("E46").Value = "  +3.33%  "
$ws.Range("D47").Value = "98.77"
$ws.Range("E47").Value = "  +0.43%  "
$ws.Range("E48").Value = "  +12.93%  "
$ws.Range("E49").Value = "  +5.60%  "
$ws.Range("D50").Value = "1.442.50"
$ws.Range("E50").Value = "  +0.62%  "
$ws.Range("D51").Value = "53.03"
$ws.Range("E51").Value = "  +7.46%  "
